$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBS")

# Append three new reference rows (FunNm / 讀取Key條件 / 其他ORDER條件)
# describing the new "Dtl" (detail) entry-date lookup functions, following the
# same pattern as the existing rows in this sheet.
$ws.Range("A5").Value = "findDtlEntryDate"
$ws.Range("B5").Value = "DtlEntryDate >= ,AND DtlEntryDate <= "
$ws.Range("C5").Value = "DtlEntryDate,Factor,CustNo,DtlSeq ASC"

$ws.Range("A6").Value = "DtlEntryDateFactor"
$ws.Range("B6").Value = "DtlEntryDate >= ,AND DtlEntryDate <= ,AND Factor ="
$ws.Range("C6").Value = "DtlEntryDate,CustNo,DtlSeq ASC"

$ws.Range("A7").Value = "findDtlEntryDateRangeFactorCustNo"
$ws.Range("B7").Value = "DtlEntryDate >= ,AND DtlEntryDate <= ,AND Factor = , AND CustNo = "
$ws.Range("C7").Value = "DtlEntryDate,DtlSeq ASC"

# Column B of this table wraps its text (matching the existing rows above).
$ws.Range("B5:B7").WrapText = $true

# Widen column A (auto-fit) so the longer new function names fit, and leave
# the cursor positioned below the newly-entered data, as Excel would after
# manual entry.
$ws.Columns.Item(1).ColumnWidth = 32.78571428571428
$ws.Range("A8").Select()
